$wb = $excel.ActiveWorkbook

# --- Sheet 2: "Displacement Table" - extend draft range from 8m to 14m ---
$wsDisp = $wb.Worksheets.Item("Displacement Table")

$dispData = @(
    @(8.5, 9800),
    @(9, 10480),
    @(9.5, 11180),
    @(10, 11900),
    @(10.5, 12640),
    @(11, 13400),
    @(11.5, 14180),
    @(12, 14980),
    @(12.5, 15800),
    @(13, 16640),
    @(13.5, 17500),
    @(14, 18380)
)

$startRow = 15
for ($i = 0; $i -lt $dispData.Count; $i++) {
    $row = $startRow + $i
    $wsDisp.Cells.Item($row, 1).Value = $dispData[$i][0]
    $wsDisp.Cells.Item($row, 2).Value = $dispData[$i][1]
}

# --- Sheet 3: "KN Curves" - add corresponding KN data rows 13-33 ---
$wsKN = $wb.Worksheets.Item("KN Curves")

$knData = @(
    @(8500, 0.78, 1.4, 2.02, 2.59, 3.11, 3.58, 4, 4.37, 4.69, 4.96, 5.1),
    @(9000, 0.81, 1.45, 2.09, 2.68, 3.22, 3.71, 4.15, 4.54, 4.88, 5.17, 5.32),
    @(9500, 0.84, 1.5, 2.16, 2.77, 3.33, 3.84, 4.3, 4.71, 5.07, 5.38, 5.54),
    @(10000, 0.87, 1.55, 2.23, 2.86, 3.44, 3.97, 4.45, 4.88, 5.26, 5.59, 5.76),
    @(10500, 0.9, 1.6, 2.3, 2.95, 3.55, 4.1, 4.6, 5.05, 5.45, 5.8, 5.98),
    @(11000, 0.93, 1.65, 2.37, 3.04, 3.66, 4.23, 4.75, 5.22, 5.64, 6.01, 6.2),
    @(11500, 0.96, 1.7, 2.44, 3.13, 3.77, 4.36, 4.9, 5.39, 5.83, 6.22, 6.42),
    @(12000, 0.99, 1.75, 2.51, 3.22, 3.88, 4.49, 5.05, 5.56, 6.02, 6.43, 6.64),
    @(12500, 1.02, 1.8, 2.58, 3.31, 3.99, 4.62, 5.2, 5.73, 6.21, 6.64, 6.86),
    @(13000, 1.05, 1.85, 2.65, 3.4, 4.1, 4.75, 5.35, 5.9, 6.4, 6.85, 7.08),
    @(13500, 1.08, 1.9, 2.72, 3.49, 4.21, 4.88, 5.5, 6.07, 6.59, 7.06, 7.3),
    @(14000, 1.11, 1.95, 2.79, 3.58, 4.32, 5.01, 5.65, 6.24, 6.78, 7.27, 7.52),
    @(14500, 1.14, 2, 2.86, 3.67, 4.43, 5.14, 5.8, 6.41, 6.97, 7.48, 7.74),
    @(15000, 1.17, 2.05, 2.93, 3.76, 4.54, 5.27, 5.95, 6.58, 7.16, 7.69, 7.96),
    @(15500, 1.2, 2.1, 3, 3.85, 4.65, 5.4, 6.1, 6.75, 7.35, 7.9, 8.18),
    @(16000, 1.23, 2.15, 3.07, 3.94, 4.76, 5.53, 6.25, 6.92, 7.54, 8.11, 8.4),
    @(16500, 1.26, 2.2, 3.14, 4.03, 4.87, 5.66, 6.4, 7.09, 7.73, 8.32, 8.62),
    @(17000, 1.29, 2.25, 3.21, 4.12, 4.98, 5.79, 6.55, 7.26, 7.92, 8.53, 8.84),
    @(17500, 1.32, 2.3, 3.28, 4.21, 5.09, 5.92, 6.7, 7.43, 8.11, 8.74, 9.06),
    @(18000, 1.35, 2.35, 3.35, 4.3, 5.2, 6.05, 6.85, 7.6, 8.3, 8.95, 9.28),
    @(18500, 1.38, 2.4, 3.42, 4.39, 5.31, 6.18, 7, 7.77, 8.49, 9.16, 9.5)
)

$knStartRow = 13
for ($i = 0; $i -lt $knData.Count; $i++) {
    $row = $knStartRow + $i
    for ($c = 0; $c -lt $knData[$i].Count; $c++) {
        $wsKN.Cells.Item($row, $c + 1).Value = $knData[$i][$c]
    }
}
